$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 1.016805354127222
$ws.Range("R2").Value = 9.151248187145001
$ws.Range("S2").Value = 0.01647576914027086
$ws.Range("T2").Value = 0.01647576914027086
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 15.16272007294089
$ws.Range("R3").Value = 136.464480656468
$ws.Range("S3").Value = 0.2456885916722543
$ws.Range("T3").Value = 0.2456885916722542
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3732763333333333
$ws.Range("H4").Value = 1.119829
$ws.Range("I4").Value = 0.3554258969843855
$ws.Range("J4").Value = 0.3554258969843855
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 5.755654167421667
$ws.Range("R4").Value = 51.800887506795
$ws.Range("S4").Value = 0.09326153617186038
$ws.Range("T4").Value = 0.09326153617186038
$ws.Range("G5").Value = 0.668317
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 1.820496621861667
$ws.Range("R5").Value = 16.384469596755
$ws.Range("S5").Value = 0.0294983518140316
$ws.Range("T5").Value = 0.0294983518140316
$ws.Range("G6").Value = 0.668317
$ws.Range("I6").Value = 0.6363574327729865
$ws.Range("J6").Value = 0.6363574327729865
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("R6").Value = 244.3271222272921
$ws.Range("S6").Value = 0.4398828638674993
$ws.Range("T6").Value = 0.4398828638674993
$ws.Range("G7").Value = 0.668317
$ws.Range("I7").Value = 0.6363574327729865
$ws.Range("J7").Value = 0.6363574327729865
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 10.304970293345
$ws.Range("R7").Value = 92.744732640105
$ws.Range("S7").Value = 0.1669762170914556
$ws.Range("T7").Value = 0.1669762170914556
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008629333333333334
$ws.Range("H8").Value = 0.025888
$ws.Range("I8").Value = 0.008216670242627913
$ws.Range("J8").Value = 0.008216670242627911
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 0.02350631838222222
$ws.Range("R8").Value = 0.21155686544
$ws.Range("S8").Value = 0.000380883788063474
$ws.Range("T8").Value = 0.0003808837880634738
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008629333333333334
$ws.Range("H9").Value = 0.025888
$ws.Range("I9").Value = 0.008216670242627913
$ws.Range("J9").Value = 0.008216670242627911
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 0.350528962232889
$ws.Range("R9").Value = 3.154760660096001
$ws.Range("S9").Value = 0.005679783485881612
$ws.Range("T9").Value = 0.00567978348588161
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008629333333333334
$ws.Range("H10").Value = 0.025888
$ws.Range("I10").Value = 0.008216670242627913
$ws.Range("J10").Value = 0.008216670242627911
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 0.1330581500266667
$ws.Range("R10").Value = 1.19752335024
$ws.Range("S10").Value = 0.002156002968682827
$ws.Range("T10").Value = 0.002156002968682827